$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 297.44446
$ws.Range("J55").Value = 314
$ws.Range("L55").Value = 314
$ws.Range("N55").Value = -742

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H80").Value = 2500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 7500
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -9496

$ws.Range("H83").Value = 2500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 22500
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -32484

$ws.Range("H92").Value = 1760.091
$ws.Range("I92").Value = 1859.125
$ws.Range("J92").Value = 1496
$ws.Range("K92").Value = 1859.125
$ws.Range("L92").Value = 1496
$ws.Range("M92").Value = -611.125
$ws.Range("N92").Value = -3992

$ws.Range("I113").Value = 55559184
$ws.Range("J113").Value = 5304.2
$ws.Range("K113").Value = 55559184
$ws.Range("L113").Value = 5304.2
$ws.Range("M113").Value = -55555930
$ws.Range("N113").Value = -11812.2

$ws.Range("H129").Value = 1227.375
$ws.Range("J129").Value = 1956.25
$ws.Range("L129").Value = 5868.75
$ws.Range("N129").Value = -15868.75

$ws.Range("H137").Value = 7454.375
$ws.Range("I137").Value = 2411.9167
$ws.Range("J137").Value = 12496.833
$ws.Range("K137").Value = 7235.750100000001
$ws.Range("L137").Value = 37490.499
$ws.Range("M137").Value = -4685.750100000001
$ws.Range("N137").Value = -42590.499

$ws.Range("H141").Value = 2790.2354
$ws.Range("I141").Value = 2790.2354
$ws.Range("K141").Value = 8370.706200000001
$ws.Range("M141").Value = -3190.706200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3252.5667
$ws.Range("I32").Value = 3261.3447
$ws.Range("K32").Value = 3261.3447
$ws.Range("M32").Value = -2974.3447

$ws.Range("H97").Value = 1289.8572
$ws.Range("I97").Value = 1296.5
$ws.Range("J97").Value = 1250
$ws.Range("K97").Value = 1296.5
$ws.Range("L97").Value = 1250
$ws.Range("M97").Value = -800.5
$ws.Range("N97").Value = -2242

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 83333620
$ws.Range("I80").Value = 333333500
$ws.Range("J80").Value = 324.1111
$ws.Range("K80").Value = 333333500
$ws.Range("L80").Value = 324.1111
$ws.Range("M80").Value = -333332502
$ws.Range("N80").Value = -2320.1111

$ws.Range("H83").Value = 83333620
$ws.Range("I83").Value = 333333500
$ws.Range("J83").Value = 324.1111
$ws.Range("K83").Value = 1666667500
$ws.Range("L83").Value = 1620.5555
$ws.Range("M83").Value = -1666662508
$ws.Range("N83").Value = -11604.5555

$ws.Range("H86").Value = 2113.4688
$ws.Range("I86").Value = 2007.1538
$ws.Range("K86").Value = 2007.1538
$ws.Range("M86").Value = -884.1538

$ws.Range("H89").Value = 2113.4688
$ws.Range("I89").Value = 2007.1538
$ws.Range("K89").Value = 10035.769
$ws.Range("M89").Value = -4419.769

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3911.4773
$ws.Range("I31").Value = 2722.5925
$ws.Range("K31").Value = 2722.5925
$ws.Range("M31").Value = -2427.5925

$ws.Range("H34").Value = 3911.4773
$ws.Range("I34").Value = 2722.5925
$ws.Range("K34").Value = 2722.5925
$ws.Range("M34").Value = -2520.5925

$ws.Range("H58").Value = 2243.0588
$ws.Range("I58").Value = 1593.3334
$ws.Range("K58").Value = 1593.3334
$ws.Range("M58").Value = -1390.3334

$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H80").Value = 49999
$ws.Range("J80").Value = 49999
$ws.Range("L80").Value = 49999
$ws.Range("N80").Value = -52245

$ws.Range("H83").Value = 49999
$ws.Range("J83").Value = 49999
$ws.Range("L83").Value = 149997
$ws.Range("N83").Value = -161229

$ws.Range("H97").Value = 44994.5
$ws.Range("J97").Value = 44994.5
$ws.Range("L97").Value = 44994.5
$ws.Range("N97").Value = -46976.5

$ws.Range("H102").Value = 26294
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 26294
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 26294
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -31162

$ws.Range("H105").Value = 2782.5715
$ws.Range("I105").Value = 2911.1667
$ws.Range("K105").Value = 2911.1667
$ws.Range("M105").Value = -1164.1667

$ws.Range("H136").Value = 2243.0588
$ws.Range("I136").Value = 1593.3334
$ws.Range("K136").Value = 4780.0002
$ws.Range("M136").Value = -2230.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 79.48484999999999
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws.Range("H5").Value = 709.5
$ws.Range("I5").Value = 709.5
$ws.Range("K5").Value = 2128.5
$ws.Range("M5").Value = -2016.5

$ws.Range("H12").Value = 5501
$ws.Range("J12").Value = 5501
$ws.Range("L12").Value = 16503
$ws.Range("N12").Value = -16849

$ws.Range("H86").Value = 1167.909
$ws.Range("I86").Value = 749.25
$ws.Range("J86").Value = 1407.1428
$ws.Range("K86").Value = 2247.75
$ws.Range("L86").Value = 4221.428400000001
$ws.Range("M86").Value = -1061.75
$ws.Range("N86").Value = -6593.428400000001

$ws.Range("H89").Value = 1167.909
$ws.Range("I89").Value = 749.25
$ws.Range("J89").Value = 1407.1428
$ws.Range("K89").Value = 6743.25
$ws.Range("L89").Value = 12664.2852
$ws.Range("M89").Value = -815.25
$ws.Range("N89").Value = -24520.2852

$ws.Range("H107").Value = 663.8421
$ws.Range("J107").Value = 830
$ws.Range("L107").Value = 2490
$ws.Range("N107").Value = -6330

$ws.Range("H132").Value = 1093.25
$ws.Range("J132").Value = 999.5
$ws.Range("L132").Value = 8995.5
$ws.Range("N132").Value = -14055.5

$ws.Range("H135").Value = 709.5
$ws.Range("I135").Value = 709.5
$ws.Range("K135").Value = 6385.5
$ws.Range("M135").Value = -3850.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5673.8
$ws.Range("I70").Value = 3532.5
$ws.Range("J70").Value = 6209.125
$ws.Range("K70").Value = 3532.5
$ws.Range("L70").Value = 6209.125
$ws.Range("M70").Value = -3262.5
$ws.Range("N70").Value = -6749.125

$ws.Range("H73").Value = 5673.8
$ws.Range("I73").Value = 3532.5
$ws.Range("J73").Value = 6209.125
$ws.Range("K73").Value = 3532.5
$ws.Range("L73").Value = 6209.125
$ws.Range("M73").Value = -2596.5
$ws.Range("N73").Value = -8081.125

$ws.Range("H97").Value = 1239.5
$ws.Range("I97").Value = 1337.8889
$ws.Range("K97").Value = 1337.8889
$ws.Range("M97").Value = -841.8888999999999

$ws.Range("H126").Value = 7921.1875
$ws.Range("J126").Value = 12154.556
$ws.Range("L126").Value = 36463.66800000001
$ws.Range("N126").Value = -41403.66800000001

$ws.Range("H134").Value = 22500
$ws.Range("J134").Value = 22500
$ws.Range("L134").Value = 67500
$ws.Range("N134").Value = -72570

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 47999
$ws.Range("J64").Value = 47999
$ws.Range("L64").Value = 47999
$ws.Range("N64").Value = -48449

$ws.Range("H67").Value = 47999
$ws.Range("J67").Value = 47999
$ws.Range("L67").Value = 47999
$ws.Range("N67").Value = -49559

$ws.Range("H80").Value = 47498.5
$ws.Range("J80").Value = 49999
$ws.Range("L80").Value = 49999
$ws.Range("N80").Value = -52245

$ws.Range("H83").Value = 47498.5
$ws.Range("J83").Value = 49999
$ws.Range("L83").Value = 149997
$ws.Range("N83").Value = -161229

$ws.Range("H93").Value = 947.5
$ws.Range("I93").Value = 395
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 395
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 853
$ws.Range("N93").Value = -3996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value = 26664.666
$ws.Range("I93").Value = 14997.5
$ws.Range("K93").Value = 14997.5
$ws.Range("M93").Value = -12501.5

$ws.Range("H96").Value = 1407.8
$ws.Range("J96").Value = 1046.3334
$ws.Range("L96").Value = 1046.3334
$ws.Range("N96").Value = -3792.3334

$ws.Range("H100").Value = 166667980
$ws.Range("I100").Value = 1603.75
$ws.Range("J100").Value = 500000740
$ws.Range("K100").Value = 3207.5
$ws.Range("L100").Value = 1000001480
$ws.Range("M100").Value = -2666.5
$ws.Range("N100").Value = -1000002562

$ws.Range("H137").Value = 77249.164
$ws.Range("J137").Value = 77249.164
$ws.Range("L137").Value = 77249.164
$ws.Range("N137").Value = -87449.164
